$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 77422
$ws.Range("B2").Value = "Dr. Bryan Moraes"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45079
$ws.Range("G2").Value = 9080.17

# Row 3
$ws.Range("A3").Value = 8255
$ws.Range("B3").Value = "Leandro Costa"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45083
$ws.Range("G3").Value = 2668.62

# Row 4
$ws.Range("A4").Value = 73494
$ws.Range("B4").Value = "Pedro Miguel Cardoso"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45103
$ws.Range("G4").Value = 7201.3

# Row 5
$ws.Range("A5").Value = 16502
$ws.Range("B5").Value = "Anthony Cavalcanti"
$ws.Range("C5").Value = "TI"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 45079
$ws.Range("G5").Value = 4261.57

# Row 6
$ws.Range("A6").Value = 29062
$ws.Range("B6").Value = "Kevin Souza"
$ws.Range("C6").Value = "Marketing"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45089
$ws.Range("G6").Value = 8706.870000000001

# Row 7
$ws.Range("A7").Value = 82467
$ws.Range("B7").Value = "Dra. Emilly Gomes"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45102
$ws.Range("G7").Value = 6898.44

# Row 8
$ws.Range("A8").Value = 31561
$ws.Range("B8").Value = "Enrico da Conceição"
$ws.Range("C8").Value = "Vendas"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45086
$ws.Range("G8").Value = 5137.35

# Row 9
$ws.Range("A9").Value = 88036
$ws.Range("B9").Value = "Sr. Benício Porto"
$ws.Range("C9").Value = "Marketing"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45095
$ws.Range("G9").Value = 4946.67

# Row 10
$ws.Range("A10").Value = 39277
$ws.Range("B10").Value = "Emilly Ferreira"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45103
$ws.Range("G10").Value = 6489.47

# Row 11
$ws.Range("A11").Value = 74754
$ws.Range("B11").Value = "Maria Alice Viana"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45078
$ws.Range("G11").Value = 3849.8

$wb.Save()
